$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.077.42"
$ws.Range("E2").Value = "  +1.61%  "
$ws.Range("D3").Value = "2.556.74"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").Value = "96.75"
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.543"
$ws.Range("E9").Value = "  +3.42%  "
$ws.Range("D10").Value = "35.62"
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("D11").Value = "0.0812"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "7.47"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("E13").Value = "  -4.77%  "
$ws.Range("D14").Value = "2.948.16"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").Value = "2.519.71"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").Value = "15.02"
$ws.Range("E16").Value = "  -1.95%  "
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "43.110.88"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").Value = "6.87"
$ws.Range("E19").Value = "  +5.15%  "
$ws.Range("D20").Value = "12.59"
$ws.Range("E20").Value = "  -3.28%  "
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D22").Value = "69.93"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "254.36"
$ws.Range("E23").Value = "  +1.57%  "
$ws.Range("D24").Value = "2.95"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  +3.19%  "
$ws.Range("D26").Value = "26.68"
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").Value = "2.44"
$ws.Range("E28").Value = "  +2.60%  "
$ws.Range("D29").Value = "40.28"
$ws.Range("E29").Value = "  +3.96%  "
$ws.Range("D30").Value = "10.25"
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("D32").Value = "155.03"
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("B33").Value = "Celestia"
$ws.Range("C33").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D33").Value = "19.23"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "2.70"
$ws.Range("E34").Value = "  +2.89%  "
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("D36").Value = "0.0803"
$ws.Range("E36").Value = "  +2.75%  "
$ws.Range("D37").Value = "3.33"
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("E38").Value = "  +1.89%  "
$ws.Range("E39").Value = "  +5.21%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "22.31"
$ws.Range("E41").Value = "  -5.22%  "
$ws.Range("D42").Value = "3.86"
$ws.Range("E42").Value = "  +2.54%  "
$ws.Range("D43").Value = "0.0304"
$ws.Range("E43").Value = "  +1.91%  "
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("D46").Value = "1.981.17"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").Value = "84.55"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").Value = "9.01"
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("D49").Value = "2.802.37"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("D50").Value = "104.33"
$ws.Range("E50").Value = "  +2.62%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.194"
$ws.Range("E51").Value = "  +3.35%  "
